# Applies the "Test for credentials AG" soapui test-results refresh:
# updates execution timestamps/durations/results on the TESTS_WS and
# TESTS_JMS sheets to reflect the latest automated test run.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TESTS_WS")
$ws2 = $wb.Worksheets.Item("TESTS_JMS")

# --- TESTS_WS sheet updates ---
$ws1.Range("G2").Value = 'PASS'
$ws1.Range("H2").Value = 42877.38479276621
$ws1.Range("I2").Value = '0.741s'
$ws1.Range("K2").Value = ''
$ws1.Range("H3").Value = 42877.384805092595
$ws1.Range("I3").Value = '9.144s'
$ws1.Range("H4").Value = 42877.384918564814
$ws1.Range("I4").Value = '2.232s'
$ws1.Range("H8").Value = 42877.3849509838
$ws1.Range("I8").Value = '12.464s'
$ws1.Range("H10").Value = 42877.385100011576
$ws1.Range("I10").Value = '2.247s'
$ws1.Range("H12").Value = 42877.38512813657
$ws1.Range("I12").Value = '6.397s'
$ws1.Range("H13").Value = 42877.385206770836
$ws1.Range("I13").Value = '4.102s'
$ws1.Range("H14").Value = 42877.385257141206
$ws1.Range("I14").Value = '2.255s'
$ws1.Range("H16").Value = 42877.38528768519
$ws1.Range("I16").Value = '4.085s'
$ws1.Range("H17").Value = 42877.38533695602
$ws1.Range("I17").Value = '4.097s'
$ws1.Range("H18").Value = 42877.3853865162
$ws1.Range("I18").Value = '4.129s'
$ws1.Range("H21").Value = 42877.3854366088
$ws1.Range("I21").Value = '4.09s'
$ws1.Range("H26").Value = 42877.38548612269
$ws1.Range("I26").Value = '4.108s'
$ws1.Range("H28").Value = 42877.38553712963
$ws1.Range("I28").Value = '4.104s'
$ws1.Range("H32").Value = 42877.38558633102
$ws1.Range("I32").Value = '4.117s'
$ws1.Range("H33").Value = 42877.385635844905
$ws1.Range("I33").Value = '4.087s'
$ws1.Range("H36").Value = 42877.385685092595
$ws1.Range("I36").Value = '4.058s'
$ws1.Range("H43").Value = 42877.38573361111
$ws1.Range("I43").Value = '4.076s'
$ws1.Range("H47").Value = 42877.385782592595
$ws1.Range("I47").Value = '4.163s'
$ws1.Range("H48").Value = 42877.38583287037
$ws1.Range("I48").Value = '4.204s'
$ws1.Range("H51").Value = 42877.38588541667
$ws1.Range("I51").Value = '4.095s'
$ws1.Range("H54").Value = 42877.38593476852
$ws1.Range("I54").Value = '4.071s'
$ws1.Range("H56").Value = 42877.385984016204
$ws1.Range("I56").Value = '4.114s'
$ws1.Range("H58").Value = 42877.38603420139
$ws1.Range("I58").Value = '4.09s'
$ws1.Range("H59").Value = 42877.38608390046
$ws1.Range("I59").Value = '4.081s'
$ws1.Range("H61").Value = 42877.38613329861
$ws1.Range("I61").Value = '2.192s'
$ws1.Range("H64").Value = 42877.38616327546
$ws1.Range("I64").Value = '4.128s'
$ws1.Range("H66").Value = 42877.38621416667
$ws1.Range("I66").Value = '4.154s'
$ws1.Range("H67").Value = 42877.386264594905
$ws1.Range("I67").Value = '4.083s'
$ws1.Range("H68").Value = 42877.3863134375
$ws1.Range("I68").Value = '4.07s'
$ws1.Range("H70").Value = 42877.386362199075
$ws1.Range("I70").Value = '4.076s'
$ws1.Range("H72").Value = 42877.38641104167
$ws1.Range("I72").Value = '4.085s'
$ws1.Range("H73").Value = 42877.38646195602
$ws1.Range("I73").Value = '4.075s'
$ws1.Range("H74").Value = 42877.38651131945
$ws1.Range("I74").Value = '4.079s'
$ws1.Range("H76").Value = 42877.386560266204
$ws1.Range("I76").Value = '4.075s'
$ws1.Range("H77").Value = 42877.38660943287
$ws1.Range("I77").Value = '4.069s'
$ws1.Range("H78").Value = 42877.386658159725
$ws1.Range("H79").Value = 42877.386706840276
$ws1.Range("I79").Value = '4.063s'
$ws1.Range("H81").Value = 42877.386755451385
$ws1.Range("H82").Value = 42877.38680457176
$ws1.Range("I82").Value = '4.111s'
$ws1.Range("H83").Value = 42877.38685380787
$ws1.Range("I83").Value = '4.067s'
$ws1.Range("H85").Value = 42877.386902673614
$ws1.Range("I85").Value = '4.1s'
$ws1.Range("H88").Value = 42877.38695216435
$ws1.Range("I88").Value = '4.089s'
$ws1.Range("H89").Value = 42877.38700123843
$ws1.Range("I89").Value = '4.082s'
$ws1.Range("H91").Value = 42877.387050300924
$ws1.Range("I91").Value = '5.491s'
$ws1.Range("H92").Value = 42877.38711710648
$ws1.Range("I92").Value = '4.342s'
$ws1.Range("H93").Value = 42877.387171041664
$ws1.Range("I93").Value = '0.025s'
$ws1.Range("H95").Value = 42877.38717255787
$ws1.Range("I95").Value = '3.429s'
$ws1.Range("H96").Value = 42877.3872147338
$ws1.Range("I96").Value = '4.165s'
$ws1.Range("H100").Value = 42877.38726611111
$ws1.Range("I100").Value = '0.023s'
$ws1.Range("H101").Value = 42877.38726761574
$ws1.Range("I101").Value = '18.166s'
$ws1.Range("H104").Value = 42877.38748409722
$ws1.Range("I104").Value = '5.172s'
$ws1.Range("H105").Value = 42877.387547025464
$ws1.Range("I105").Value = '4.06s'
$ws1.Range("H107").Value = 42877.38759568287
$ws1.Range("I107").Value = '2.19s'
$ws1.Range("H108").Value = 42877.38762424768
$ws1.Range("I108").Value = '2.119s'
$ws1.Range("H110").Value = 42877.38765089121
$ws1.Range("I110").Value = '8.087s'
$ws1.Range("H113").Value = 42877.387747534725
$ws1.Range("I113").Value = '0.095s'
$ws1.Range("H118").Value = 42877.387749988426
$ws1.Range("I118").Value = '4.238s'

# --- TESTS_JMS sheet updates ---
$ws2.Range("H2").Value = 42881.558739618056
$ws2.Range("I2").Value = '0.181s'
$ws2.Range("H3").Value = 42881.55874328704
$ws2.Range("I3").Value = '0.069s'
$ws2.Range("H4").Value = 42881.55874626157
$ws2.Range("I4").Value = '2.168s'
$ws2.Range("I5").Value = '14.131s'
$ws2.Range("H6").Value = 42881.56552023148
$ws2.Range("I6").Value = '195.108s'
$ws2.Range("H8").Value = 42881.55877556713
$ws2.Range("I8").Value = '6.282s'
$ws2.Range("H9").Value = 42881.55885269676
$ws2.Range("I9").Value = '4.116s'
$ws2.Range("H10").Value = 42881.55890246528
$ws2.Range("I10").Value = '4.224s'
$ws2.Range("H11").Value = 42881.55895494213
$ws2.Range("I11").Value = '4.139s'
$ws2.Range("H14").Value = 42881.55900576389
$ws2.Range("I14").Value = '4.14s'
$ws2.Range("H17").Value = 42881.5590553125
$ws2.Range("I17").Value = '4.146s'
$ws2.Range("H19").Value = 42881.55911054398
$ws2.Range("I19").Value = '13.062s'
$ws2.Range("H20").Value = 42881.55926628472
$ws2.Range("I20").Value = '12.758s'
$ws2.Range("H26").Value = 42881.5594183912
$ws2.Range("I26").Value = '8.284s'
$ws2.Range("H27").Value = 42881.55951856481
$ws2.Range("I27").Value = '4.162s'
$ws2.Range("H31").Value = 42881.559568715275
$ws2.Range("I31").Value = '4.507s'
$ws2.Range("H34").Value = 42881.55962290509
$ws2.Range("I34").Value = '8.244s'
$ws2.Range("H36").Value = 42881.55972208334
$ws2.Range("I36").Value = '8.28s'
$ws2.Range("H38").Value = 42881.559821875
$ws2.Range("I38").Value = '8.293s'
$ws2.Range("H40").Value = 42881.559922280096
$ws2.Range("I40").Value = '8.237s'
$ws2.Range("H46").Value = 42881.560020763885
$ws2.Range("I46").Value = '8.22s'
$ws2.Range("D50").Value = 'DomJMS048-Submit Message-payload_contentId'
$ws2.Range("H50").Value = 42881.560122118055
$ws2.Range("I50").Value = '12.632s'
$ws2.Range("G52").Value = 'FAIL'
$ws2.Range("H52").Value = 42881.56173534722
$ws2.Range("I52").Value = '2.133s'
$ws2.Range("K52").Value = '26-05-2017 13:28:54: Test case FAILED on step 2: Test|| Returned error message[s]: 
 |java.lang.AssertionError: --DomJMS050-Test Experience-Several payloads--Test--  Error:verifyMessagePresence: Message with ID -DomJMS050- is not found in sender side.. Expression: (total > 0). Values: total = 0| 
 |error at line: 2| '
$ws2.Range("H53").Value = 42881.56029625
$ws2.Range("I53").Value = '4.096s'
$ws2.Range("H54").Value = 42881.56947436342
$ws2.Range("I54").Value = '55.463s'
$ws2.Range("H5").Value = 42881.564579224534

